$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item("Subtitle 2")
$tf = $sh.TextFrame
$tr = $tf.TextRange
$para2 = $tr.Paragraphs(2)
$r1 = $para2.Runs(1)
$r1.Text = "نام "
$r1.InsertAfter("استاد:دکتر عصایی") | Out-Null
